{"js": "// Replace the twenty-five \"two-digit \u00d7 two-digit = product\" answer cells\n// in the multiplication worksheet table with a new set of equations,\n// per the commit diff. Each old value is unique in the document, so a\n// plain text search-and-replace (scoped to exact, case-sensitive matches)\n// is unambiguous.\nconst replacements = [\n  [\"27\u00d765=1755\", \"98\u00d714=1372\"],\n  [\"72\u00d739=2808\", \"50\u00d749=2450\"],\n  [\"14\u00d711=154\", \"71\u00d769=4899\"],\n  [\"54\u00d790=4860\", \"45\u00d799=4455\"],\n  [\"60\u00d755=3300\", \"46\u00d784=3864\"],\n  [\"46\u00d731=1426\", \"78\u00d763=4914\"],\n  [\"63\u00d741=2583\", \"45\u00d795=4275\"],\n  [\"78\u00d774=5772\", \"80\u00d715=1200\"],\n  [\"64\u00d746=2944\", \"76\u00d730=2280\"],\n  [\"73\u00d797=7081\", \"74\u00d736=2664\"],\n  [\"51\u00d744=2244\", \"19\u00d774=1406\"],\n  [\"24\u00d737=888\", \"49\u00d753=2597\"],\n  [\"46\u00d717=782\", \"34\u00d738=1292\"],\n  [\"57\u00d783=4731\", \"14\u00d741=574\"],\n  [\"12\u00d769=828\", \"52\u00d778=4056\"],\n  [\"22\u00d770=1540\", \"72\u00d795=6840\"],\n  [\"15\u00d787=1305\", \"96\u00d780=7680\"],\n  [\"23\u00d770=1610\", \"55\u00d799=5445\"],\n  [\"73\u00d721=1533\", \"28\u00d772=2016\"],\n  [\"33\u00d711=363\", \"60\u00d722=1320\"],\n  [\"88\u00d756=4928\", \"47\u00d725=1175\"],\n  [\"15\u00d761=915\", \"99\u00d795=9405\"],\n  [\"54\u00d764=3456\", \"69\u00d781=5589\"],\n  [\"94\u00d754=5076\", \"40\u00d739=1560\"],\n  [\"70\u00d783=5810\", \"16\u00d726=416\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the twenty-five \"two-digit \u00d7 two-digit = product\" answer cells\n# in the multiplication worksheet table with a new set of equations, per\n# the commit diff. Each old value is unique in the document, so a plain\n# Find/Replace (MatchCase, no wildcards) for each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"27\u00d765=1755\", \"98\u00d714=1372\"),\n  @(\"72\u00d739=2808\", \"50\u00d749=2450\"),\n  @(\"14\u00d711=154\",  \"71\u00d769=4899\"),\n  @(\"54\u00d790=4860\", \"45\u00d799=4455\"),\n  @(\"60\u00d755=3300\", \"46\u00d784=3864\"),\n  @(\"46\u00d731=1426\", \"78\u00d763=4914\"),\n  @(\"63\u00d741=2583\", \"45\u00d795=4275\"),\n  @(\"78\u00d774=5772\", \"80\u00d715=1200\"),\n  @(\"64\u00d746=2944\", \"76\u00d730=2280\"),\n  @(\"73\u00d797=7081\", \"74\u00d736=2664\"),\n  @(\"51\u00d744=2244\", \"19\u00d774=1406\"),\n  @(\"24\u00d737=888\",  \"49\u00d753=2597\"),\n  @(\"46\u00d717=782\",  \"34\u00d738=1292\"),\n  @(\"57\u00d783=4731\", \"14\u00d741=574\"),\n  @(\"12\u00d769=828\",  \"52\u00d778=4056\"),\n  @(\"22\u00d770=1540\", \"72\u00d795=6840\"),\n  @(\"15\u00d787=1305\", \"96\u00d780=7680\"),\n  @(\"23\u00d770=1610\", \"55\u00d799=5445\"),\n  @(\"73\u00d721=1533\", \"28\u00d772=2016\"),\n  @(\"33\u00d711=363\",  \"60\u00d722=1320\"),\n  @(\"88\u00d756=4928\", \"47\u00d725=1175\"),\n  @(\"15\u00d761=915\",  \"99\u00d795=9405\"),\n  @(\"54\u00d764=3456\", \"69\u00d781=5589\"),\n  @(\"94\u00d754=5076\", \"40\u00d739=1560\"),\n  @(\"70\u00d783=5810\", \"16\u00d726=416\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
